# Auto-generated update script: applies the "Update automàtic: dades i banners [2026-02-25 07:20]" edit
# to the meteocat daily summary sheet (updates DATA_EXTRACCIO timestamps and refreshed observation values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage-looking text values (e.g. "58%") must be forced to Text format first,
# otherwise Excel auto-converts them to a numeric percentage (0.58) instead of keeping the literal string.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = '58%'
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = '20%'
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = '97%'
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = '54%'
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '94%'
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '45%'
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '78%'
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = '39%'
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = '24%'
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '91%'
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = '69%'
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = '38%'
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = '39%'
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = '93%'

# Remaining refreshed values (timestamps, pressures, wind, temperatures, snow depth, radiation, humidity).
$ws.Range("E2").Value = '2026-02-25 07:18:34'
$ws.Range("E3").Value = '2026-02-25 07:18:37'
$ws.Range("E4").Value = '2026-02-25 07:18:39'
$ws.Range("J4").Value = '1020.1 hPa'
$ws.Range("E5").Value = '2026-02-25 07:18:42'
$ws.Range("M5").Value = '6.3 °C 6:51 TU'
$ws.Range("E6").Value = '2026-02-25 07:18:45'
$ws.Range("J6").Value = '1019.6 hPa'
$ws.Range("E7").Value = '2026-02-25 07:18:47'
$ws.Range("J7").Value = '1019.0 hPa'
$ws.Range("K7").Value = '0.0 MJ/m2'
$ws.Range("L7").Value = '16.6 km/h - 34º 6:51 TU'
$ws.Range("E8").Value = '2026-02-25 07:18:50'
$ws.Range("J8").Value = '1018.4 hPa'
$ws.Range("K8").Value = '0.0 MJ/m2'
$ws.Range("N8").Value = '7.9 °C 6:41 TU'
$ws.Range("O8").Value = '13.8 °C'
$ws.Range("E9").Value = '2026-02-25 07:18:52'
$ws.Range("E10").Value = '2026-02-25 07:18:55'
$ws.Range("K10").Value = '0.0 MJ/m2'
$ws.Range("L10").Value = '6.1 km/h - 237º 6:32 TU'
$ws.Range("O10").Value = '4.7 °C'
$ws.Range("E11").Value = '2026-02-25 07:18:58'
$ws.Range("N11").Value = '1.6 °C 6:52 TU'
$ws.Range("O11").Value = '2.7 °C'
$ws.Range("E12").Value = '2026-02-25 07:19:00'
$ws.Range("M12").Value = '7.2 °C 6:41 TU'
$ws.Range("O12").Value = '5.7 °C'
$ws.Range("E13").Value = '2026-02-25 07:19:03'
$ws.Range("J13").Value = '1027.1 hPa'
$ws.Range("O13").Value = '-2.1 °C'
$ws.Range("E14").Value = '2026-02-25 07:19:05'
$ws.Range("O14").Value = '5.3 °C'
$ws.Range("E15").Value = '2026-02-25 07:19:07'
$ws.Range("N15").Value = '3.1 °C 6:47 TU'
$ws.Range("O15").Value = '5.0 °C'
$ws.Range("E16").Value = '2026-02-25 07:19:10'
$ws.Range("O16").Value = '3.5 °C'
$ws.Range("E17").Value = '2026-02-25 07:19:13'
$ws.Range("K17").Value = '0.0 MJ/m2'
$ws.Range("O17").Value = '8.7 °C'
$ws.Range("E18").Value = '2026-02-25 07:19:15'
$ws.Range("J18").Value = '1020.0 hPa'
$ws.Range("E19").Value = '2026-02-25 07:19:18'
$ws.Range("L19").Value = '11.5 km/h - 101º 6:52 TU'
$ws.Range("E20").Value = '2026-02-25 07:19:20'
$ws.Range("K20").Value = '0.0 MJ/m2'
$ws.Range("N20").Value = '0.3 °C 6:51 TU'
$ws.Range("O20").Value = '2.5 °C'
$ws.Range("E21").Value = '2026-02-25 07:19:23'
$ws.Range("J21").Value = '1024.0 hPa'
$ws.Range("K21").Value = '0.0 MJ/m2'
$ws.Range("N21").Value = '0.9 °C 6:50 TU'
$ws.Range("O21").Value = '3.0 °C'
$ws.Range("E22").Value = '2026-02-25 07:19:26'
$ws.Range("K22").Value = '0.0 MJ/m2'
$ws.Range("E23").Value = '2026-02-25 07:19:28'
$ws.Range("L23").Value = '17.6 km/h - 43º 6:33 TU'
$ws.Range("M23").Value = '4.5 °C 6:35 TU'
$ws.Range("E24").Value = '2026-02-25 07:19:31'
$ws.Range("N24").Value = '2.1 °C 6:51 TU'
$ws.Range("O24").Value = '3.5 °C'
$ws.Range("E25").Value = '2026-02-25 07:19:34'
$ws.Range("K25").Value = '0.0 MJ/m2'
$ws.Range("O25").Value = '2.9 °C'
$ws.Range("E26").Value = '2026-02-25 07:19:36'
$ws.Range("G26").Value = '1 cm'
$ws.Range("J26").Value = '1019.3 hPa'
$ws.Range("N26").Value = '7.1 °C 6:33 TU'
$ws.Range("O26").Value = '8.8 °C'
$ws.Range("E27").Value = '2026-02-25 07:19:39'
$ws.Range("E28").Value = '2026-02-25 07:19:41'
$ws.Range("J28").Value = '1020.9 hPa'
$ws.Range("N28").Value = '1.7 °C 6:33 TU'
$ws.Range("E29").Value = '2026-02-25 07:19:44'
$ws.Range("O29").Value = '9.8 °C'
$ws.Range("E30").Value = '2026-02-25 07:19:47'
$ws.Range("J30").Value = '1019.9 hPa'
$ws.Range("O30").Value = '7.2 °C'
$ws.Range("E31").Value = '2026-02-25 07:19:49'
$ws.Range("J31").Value = '1019.1 hPa'
$ws.Range("E32").Value = '2026-02-25 07:19:52'
$ws.Range("K32").Value = '0.0 MJ/m2'
$ws.Range("O32").Value = '1.8 °C'
$ws.Range("E33").Value = '2026-02-25 07:19:54'
$ws.Range("J33").Value = '1024.4 hPa'
$ws.Range("O33").Value = '1.8 °C'
$ws.Range("E34").Value = '2026-02-25 07:19:57'
$ws.Range("O34").Value = '1.4 °C'
$ws.Range("E35").Value = '2026-02-25 07:20:00'
$ws.Range("E36").Value = '2026-02-25 07:20:02'
$ws.Range("J36").Value = '1019.6 hPa'
$ws.Range("E37").Value = '2026-02-25 07:20:05'
$ws.Range("J37").Value = '1024.4 hPa'
$ws.Range("L37").Value = '10.4 km/h - 67º 6:51 TU'
$ws.Range("M37").Value = '3.9 °C 6:59 TU'
$ws.Range("E38").Value = '2026-02-25 07:20:08'
$ws.Range("O38").Value = '4.2 °C'
$ws.Range("E39").Value = '2026-02-25 07:20:10'
$ws.Range("K39").Value = '0.0 MJ/m2'
$ws.Range("E40").Value = '2026-02-25 07:20:12'
$ws.Range("M40").Value = '5.1 °C 6:38 TU'
$ws.Range("O40").Value = '1.6 °C'
$ws.Range("E41").Value = '2026-02-25 07:20:15'
$ws.Range("N41").Value = '8.4 °C 6:58 TU'
$ws.Range("E42").Value = '2026-02-25 07:20:18'
$ws.Range("E43").Value = '2026-02-25 07:20:20'
$ws.Range("K43").Value = '0.0 MJ/m2'
$ws.Range("N43").Value = '1.3 °C 6:52 TU'
$ws.Range("O43").Value = '3.3 °C'
$ws.Range("E44").Value = '2026-02-25 07:20:23'
$ws.Range("K44").Value = '0.0 MJ/m2'
$ws.Range("E45").Value = '2026-02-25 07:20:25'
$ws.Range("K45").Value = '0.0 MJ/m2'
$ws.Range("E46").Value = '2026-02-25 07:20:28'
$ws.Range("J46").Value = '1019.7 hPa'
$ws.Range("K46").Value = '0.0 MJ/m2'
$ws.Range("O46").Value = '3.1 °C'
